# Apply the "dotenv + separate body/subject file" edit to the Names sheet:
#  - Replace the sample email address and name with the new values
#  - Clear the sample name that used to sit in B3
#  - Remove the mailto: hyperlink that was attached to the email cell
#  - Mark the email cell (A2) as wrapping text, like the row below it
#  - Move/extend the active selection to A3:B3
#  - Restore header/footer margins to 1.3cm (their precise inch value)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old mailto hyperlink on A2 before changing its value.
$ws.Hyperlinks.Delete()

# Update the sample data.
$ws.Range("A2").Value = "20190016936@my.xu.edu.ph"
$ws.Range("B2").Value = "Josh Paculba"
$ws.Range("B3").Value = ""

# A2 should wrap like the cells in row 3 already do.
$ws.Range("A2").WrapText = $true

# Park the selection on the now-empty A3:B3 row.
$ws.Range("A3:B3").Select()

# Header/footer margins recompute to the precise 1.3cm-in-inches value.
$ws.PageSetup.HeaderMargin = 36.850393700787386
$ws.PageSetup.FooterMargin = 36.850393700787386
